$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# ---- Row 5 ----
$ws.Range("A5").Value = 112307840
$ws.Range("B5").Value = 89049
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5734
$ws.Range("F5").Value = "Druvfingersvamp"
$ws.Range("G5").Value = "Ramaria botrytis"
$ws.Range("H5").Value = "(Pers.:Fr.) Bourdot"
$ws.Range("P5").Value = "Hagalund, skogen norr, Adelsö, Upl"
$ws.Range("Q5").Value = 641395
$ws.Range("R5").Value = 6585171
$ws.Range("S5").Value = 25
$ws.Range("T5").Value = "Stockholm"
$ws.Range("U5").Value = "Ekerö"
$ws.Range("V5").Value = "Uppland"
$ws.Range("W5").Value = "Adelsö"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-24"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-24"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AW5").Value = "Hasse Andersson"
$ws.Range("AX5").Value = "Hasse Andersson, Birgit Svensson"
$ws.Range("Y5").ClearFormats()
$ws.Range("AA5").ClearFormats()

# ---- Row 6 ----
$ws.Range("A6").Value = 112326188
$ws.Range("B6").Value = 90800
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "LC"
$ws.Range("E6").Value = 4364
$ws.Range("F6").Value = "Dropptaggsvamp"
$ws.Range("G6").Value = "Hydnellum ferrugineum"
$ws.Range("H6").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("P6").Value = "Hagalund, skogen norr, Adelsö, Upl"
$ws.Range("Q6").Value = 641395
$ws.Range("R6").Value = 6585171
$ws.Range("S6").Value = 25
$ws.Range("T6").Value = "Stockholm"
$ws.Range("U6").Value = "Ekerö"
$ws.Range("V6").Value = "Uppland"
$ws.Range("W6").Value = "Adelsö"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-24"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-24"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AW6").Value = "Hasse Andersson"
$ws.Range("AX6").Value = "Hasse Andersson, Birgit Svensson"
$ws.Range("Y6").ClearFormats()
$ws.Range("AA6").ClearFormats()
